$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("I7").Value = "ba"
$ws.Range("J7").Value = "Appreciation"

$ws.Range("I11").Value = "aa"
$ws.Range("J11").Value = "Agree/Accept"

$ws.Range("I13").Value = "sd"
$ws.Range("J13").Value = "Statement-non-opinion"

$ws.Range("I24").Value = "sv"
$ws.Range("J24").Value = "Statement-opinion"

$ws.Range("I39").Value = "sv"
$ws.Range("J39").Value = "Statement-opinion"

$ws.Range("I75").Value = "sd"
$ws.Range("J75").Value = "Statement-non-opinion"
